# Test data for Spain added:
# duplicate the "Italy" sheet to create "Spain", then update the
# market name / ticket-reference cells on the new sheet.

$wb = $excel.ActiveWorkbook
$italy = $wb.Worksheets.Item("Italy")

# Right-click "Italy" tab -> Move or Copy -> Create a copy, placed after Italy.
$italy.Copy([System.Reflection.Missing]::Value, $italy)
$spain = $wb.Worksheets.Item("Italy (2)")
$spain.Name = "Spain"

# Update the copied sheet's market / reference cells for Spain.
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3442/T1602/T1643"

# Restore the original sheet's whole-sheet selection state, then leave
# "Spain" as the active tab with its original B4 selection.
$italy.Activate()
[void]$italy.Cells.Select()
$spain.Activate()
[void]$spain.Range("B4").Select()
